$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

# Row 2: Bitcoin -> Bitcoin
Set-TextCell $ws.Range("D2") "91.777.58"
Set-TextCell $ws.Range("E2") "  +1.19%  "

# Row 3: Ethereum -> Ethereum
Set-TextCell $ws.Range("D3") "3.129.32"
Set-TextCell $ws.Range("E3") "  -0.47%  "

# Row 4: TetherUSD -> TetherUSD
Set-TextCell $ws.Range("D4") "0.999"
Set-TextCell $ws.Range("E4") "  -0.03%  "

# Row 5: Solana -> Solana
Set-TextCell $ws.Range("D5") "240.92"
Set-TextCell $ws.Range("E5") "  +1.30%  "

# Row 6: BNB -> BNB
Set-TextCell $ws.Range("D6") "627.06"
Set-TextCell $ws.Range("E6") "  -2.24%  "

# Row 7: XRP -> XRP
Set-TextCell $ws.Range("E7") "  +9.67%  "

# Row 8: Dogecoin -> Dogecoin
Set-TextCell $ws.Range("E8") "  +3.44%  "

# Row 9: USDC -> USDC
Set-TextCell $ws.Range("E9") "  -0.03%  "

# Row 10: LidoStakedEther -> LidoStakedEther
Set-TextCell $ws.Range("D10") "3.127.40"
Set-TextCell $ws.Range("E10") "  -0.21%  "

# Row 11: Cardano -> Cardano
Set-TextCell $ws.Range("D11") "0.771"
Set-TextCell $ws.Range("E11") "  +7.10%  "

# Row 12: TRON -> TRON
Set-TextCell $ws.Range("E12") "  +3.78%  "

# Row 13: ShibaInu -> ShibaInu
Set-TextCell $ws.Range("D13") "0.0000255"
Set-TextCell $ws.Range("E13") "  +3.49%  "

# Row 14: Avalanche -> Avalanche
Set-TextCell $ws.Range("D14") "35.91"
Set-TextCell $ws.Range("E14") "  -1.77%  "

# Row 15: Toncoin -> Toncoin
Set-TextCell $ws.Range("D15") "5.53"
Set-TextCell $ws.Range("E15") "  -1.55%  "

# Row 16: WrappedBTC -> WrappedBTC
Set-TextCell $ws.Range("D16") "91.557.73"
Set-TextCell $ws.Range("E16") "  +1.24%  "

# Row 17: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
Set-TextCell $ws.Range("D17") "3.706.68"
Set-TextCell $ws.Range("E17") "  -0.20%  "

# Row 18: SuiNetwork -> WrappedEther
Set-TextCell $ws.Range("B18") "WrappedEther"
Set-TextCell $ws.Range("C18") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell $ws.Range("D18") "3.121.19"
Set-TextCell $ws.Range("E18") "  -1.42%  "

# Row 19: WrappedEther -> SuiNetwork
Set-TextCell $ws.Range("B19") "SuiNetwork"
Set-TextCell $ws.Range("C19") "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextCell $ws.Range("D19") "3.77"
Set-TextCell $ws.Range("E19") "  +0.15%  "

# Row 20: Chainlink -> Chainlink
Set-TextCell $ws.Range("D20") "14.88"
Set-TextCell $ws.Range("E20") "  +2.80%  "

# Row 21: PEPE -> PEPE
Set-TextCell $ws.Range("D21") "0.0000219"
Set-TextCell $ws.Range("E21") "  +0.64%  "

# Row 22: Polkadot -> Polkadot
Set-TextCell $ws.Range("D22") "5.93"
Set-TextCell $ws.Range("E22") "  +4.42%  "

# Row 23: BitcoinCash -> BitcoinCash
Set-TextCell $ws.Range("D23") "449.69"
Set-TextCell $ws.Range("E23") "  +0.13%  "

# Row 24: Uniswap -> Uniswap
Set-TextCell $ws.Range("D24") "9.16"
Set-TextCell $ws.Range("E24") "  +1.05%  "

# Row 25: NEARProtocol -> NEARProtocol
Set-TextCell $ws.Range("D25") "5.91"
Set-TextCell $ws.Range("E25") "  +1.53%  "

# Row 26: Litecoin -> Litecoin
Set-TextCell $ws.Range("D26") "93.25"
Set-TextCell $ws.Range("E26") "  +2.49%  "

# Row 27: Aptos -> Aptos
Set-TextCell $ws.Range("D27") "12.04"
Set-TextCell $ws.Range("E27") "  -3.00%  "

# Row 28: WrappedeETH -> WrappedeETH
Set-TextCell $ws.Range("D28") "3.297.57"
Set-TextCell $ws.Range("E28") "  +0.39%  "

# Row 29: Dai -> Dai
Set-TextCell $ws.Range("E29") "  +0.05%  "

# Row 30: Cronos -> Stellar
Set-TextCell $ws.Range("B30") "Stellar"
Set-TextCell $ws.Range("C30") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws.Range("D30") "0.248"
Set-TextCell $ws.Range("E30") "  +25.35%  "

# Row 31: Stellar -> Cronos
Set-TextCell $ws.Range("B31") "Cronos"
Set-TextCell $ws.Range("C31") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell $ws.Range("D31") "0.184"
Set-TextCell $ws.Range("E31") "  +14.92%  "

# Row 32: Hedera -> Hedera
Set-TextCell $ws.Range("E32") "  +38.83%  "

# Row 33: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
Set-TextCell $ws.Range("E33") "  -3.90%  "

# Row 34: Binance-PegBSC-USD -> Binance-PegBSC-USD
Set-TextCell $ws.Range("D34") "1.00"
Set-TextCell $ws.Range("E34") "  +23.75%  "

# Row 35: Kaspa -> Kaspa
Set-TextCell $ws.Range("D35") "0.167"
Set-TextCell $ws.Range("E35") "  +11.36%  "

# Row 36: EthereumClassic -> EthereumClassic
Set-TextCell $ws.Range("D36") "26.82"
Set-TextCell $ws.Range("E36") "  -2.03%  "

# Row 37: RenderToken -> RenderToken
Set-TextCell $ws.Range("D37") "7.57"
Set-TextCell $ws.Range("E37") "  +5.54%  "

# Row 38: MantraDAO -> MantraDAO
Set-TextCell $ws.Range("E38") "  +27.74%  "

# Row 39: Bittensor -> PancakeSwap
Set-TextCell $ws.Range("B39") "PancakeSwap"
Set-TextCell $ws.Range("C39") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws.Range("D39") "1.94"
Set-TextCell $ws.Range("E39") "  -0.96%  "

# Row 40: PancakeSwap -> Bittensor
Set-TextCell $ws.Range("B40") "Bittensor"
Set-TextCell $ws.Range("C40") "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell $ws.Range("D40") "498.37"
Set-TextCell $ws.Range("E40") "  -3.89%  "

# Row 41: dogwifhat -> dogwifhat
Set-TextCell $ws.Range("D41") "3.65"
Set-TextCell $ws.Range("E41") "  -5.25%  "

# Row 42: Fetch.AI -> Fetch.AI
Set-TextCell $ws.Range("E42") "  -0.38%  "

# Row 43: PolygonEcosystemToken -> PolygonEcosystemToken
Set-TextCell $ws.Range("D43") "0.424"
Set-TextCell $ws.Range("E43") "  +0.24%  "

# Row 44: WhiteBITCoin -> WhiteBITCoin
Set-TextCell $ws.Range("E44") "  -0.20%  "

# Row 46: Stacks -> Stacks
Set-TextCell $ws.Range("E46") "  +0.00%  "

# Row 47: Monero -> Monero
Set-TextCell $ws.Range("D47") "156.57"
Set-TextCell $ws.Range("E47") "  +3.90%  "

# Row 48: ARBITRUM -> ARBITRUM
Set-TextCell $ws.Range("D48") "0.699"
Set-TextCell $ws.Range("E48") "  -0.39%  "

# Row 49: Filecoin -> ImmutableX
Set-TextCell $ws.Range("B49") "ImmutableX"
Set-TextCell $ws.Range("C49") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell $ws.Range("D49") "1.37"
Set-TextCell $ws.Range("E49") "  +0.77%  "

# Row 50: ImmutableX -> Filecoin
Set-TextCell $ws.Range("B50") "Filecoin"
Set-TextCell $ws.Range("C50") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws.Range("D50") "4.56"
Set-TextCell $ws.Range("E50") "  +0.00%  "

# Row 51: OKB -> OKB
Set-TextCell $ws.Range("D51") "44.82"
